$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1562.0588
$ws.Range("I19").Value = 1857.2
$ws.Range("J19").Value = 1140.4286
$ws.Range("K19").Value = 1857.2
$ws.Range("L19").Value = 1140.4286
$ws.Range("M19").Value = -1682.2
$ws.Range("N19").Value = -1490.4286
$ws.Range("H76").Value = 4165.6665
$ws.Range("J76").Value = 4500
$ws.Range("L76").Value = 4500
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 4165.6665
$ws.Range("J79").Value = 4500
$ws.Range("L79").Value = 4500
$ws.Range("N79").Value = -6684
$ws.Range("H86").Value = 2768
$ws.Range("I86").Value = 2550
$ws.Range("J86").Value = 3204
$ws.Range("K86").Value = 2550
$ws.Range("L86").Value = 3204
$ws.Range("M86").Value = -1427
$ws.Range("N86").Value = -5450
$ws.Range("H89").Value = 2768
$ws.Range("I89").Value = 2550
$ws.Range("J89").Value = 3204
$ws.Range("K89").Value = 12750
$ws.Range("L89").Value = 16020
$ws.Range("M89").Value = -7134
$ws.Range("N89").Value = -27252
$ws.Range("H112").Value = 1644.3182
$ws.Range("J112").Value = 1782.8948
$ws.Range("L112").Value = 5348.6844
$ws.Range("N112").Value = -7564.6844
$ws.Range("H116").Value = 6860.25
$ws.Range("J116").Value = 6813
$ws.Range("L116").Value = 6813
$ws.Range("N116").Value = -13697
$ws.Range("H138").Value = 2759.2778
$ws.Range("I138").Value = 1834.2727
$ws.Range("J138").Value = 3166.28
$ws.Range("K138").Value = 5502.8181
$ws.Range("L138").Value = 9498.84
$ws.Range("M138").Value = -362.8181000000004
$ws.Range("N138").Value = -19778.84
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1292.75
$ws.Range("I2").Value = 1371.4
$ws.Range("J2").Value = 899.5
$ws.Range("K2").Value = 1371.4
$ws.Range("L2").Value = 899.5
$ws.Range("M2").Value = -1258.4
$ws.Range("N2").Value = -1125.5
$ws.Range("H88").Value = 1250
$ws.Range("J88").Value = 1000
$ws.Range("L88").Value = 1000
$ws.Range("N88").Value = -1812
$ws.Range("H91").Value = 1250
$ws.Range("J91").Value = 1000
$ws.Range("L91").Value = 1000
$ws.Range("N91").Value = -3808
$ws.Range("H116").Value = 1292.75
$ws.Range("I116").Value = 1371.4
$ws.Range("J116").Value = 899.5
$ws.Range("K116").Value = 1371.4
$ws.Range("L116").Value = 899.5
$ws.Range("M116").Value = 922.5999999999999
$ws.Range("N116").Value = -5487.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1292.75
$ws.Range("I3").Value = 1371.4
$ws.Range("J3").Value = 899.5
$ws.Range("K3").Value = 1371.4
$ws.Range("L3").Value = 899.5
$ws.Range("M3").Value = -1257.4
$ws.Range("N3").Value = -1127.5
$ws.Range("H20").Value = 3937.5
$ws.Range("I20").Value = 370
$ws.Range("J20").Value = 5126.6665
$ws.Range("K20").Value = 370
$ws.Range("L20").Value = 5126.6665
$ws.Range("M20").Value = -123
$ws.Range("N20").Value = -5620.6665
$ws.Range("H105").Value = 3793.5715
$ws.Range("I105").Value = 3175.8333
$ws.Range("K105").Value = 3175.8333
$ws.Range("M105").Value = -1428.8333
$ws.Range("H134").Value = 2234.8147
$ws.Range("I134").Value = 2045.0476
$ws.Range("K134").Value = 6135.142800000001
$ws.Range("M134").Value = -3600.142800000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 492
$ws.Range("I7").Value = 348.75
$ws.Range("J7").Value = 683
$ws.Range("K7").Value = 348.75
$ws.Range("L7").Value = 683
$ws.Range("M7").Value = -235.75
$ws.Range("N7").Value = -909
$ws.Range("H16").Value = 2204.6667
$ws.Range("I16").Value = 2212
$ws.Range("J16").Value = 2190
$ws.Range("K16").Value = 2212
$ws.Range("L16").Value = 2190
$ws.Range("M16").Value = -1925
$ws.Range("N16").Value = -2764
$ws.Range("H58").Value = 2723.5715
$ws.Range("I58").Value = 2523.8333
$ws.Range("J58").Value = 2873.375
$ws.Range("K58").Value = 2523.8333
$ws.Range("L58").Value = 2873.375
$ws.Range("M58").Value = -2320.8333
$ws.Range("N58").Value = -3279.375
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31498
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -97488
$ws.Range("H105").Value = 1734.375
$ws.Range("I105").Value = 1734.375
$ws.Range("K105").Value = 1734.375
$ws.Range("M105").Value = 12.625
$ws.Range("H113").Value = 2204.6667
$ws.Range("I113").Value = 2212
$ws.Range("J113").Value = 2190
$ws.Range("K113").Value = 2212
$ws.Range("L113").Value = 2190
$ws.Range("M113").Value = -42
$ws.Range("N113").Value = -6530
$ws.Range("H132").Value = 3243.25
$ws.Range("I132").Value = 2394.4443
$ws.Range("J132").Value = 4334.5713
$ws.Range("K132").Value = 7183.3329
$ws.Range("L132").Value = 13003.7139
$ws.Range("M132").Value = -4653.3329
$ws.Range("N132").Value = -18063.7139
$ws.Range("H134").Value = 4498.8
$ws.Range("I134").Value = 4498.5
$ws.Range("K134").Value = 13495.5
$ws.Range("M134").Value = -10960.5
$ws.Range("H136").Value = 2723.5715
$ws.Range("I136").Value = 2523.8333
$ws.Range("J136").Value = 2873.375
$ws.Range("K136").Value = 7571.499899999999
$ws.Range("L136").Value = 8620.125
$ws.Range("M136").Value = -5021.499899999999
$ws.Range("N136").Value = -13720.125
$ws.Range("H141").Value = 68971
$ws.Range("J141").Value = 58628
$ws.Range("L141").Value = 58628
$ws.Range("N141").Value = -68988
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1087.5
$ws.Range("J5").Value = 1350
$ws.Range("L5").Value = 4050
$ws.Range("N5").Value = -4274
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("N37").Value = 0
$ws.Range("H56").Value = 18477.824
$ws.Range("I56").Value = 18477.824
$ws.Range("K56").Value = 18477.824
$ws.Range("M56").Value = -17947.824
$ws.Range("H135").Value = 1087.5
$ws.Range("J135").Value = 1350
$ws.Range("L135").Value = 12150
$ws.Range("N135").Value = -17220
$ws.Range("H137").Value = 4179.6
$ws.Range("I137").Value = 2879.8
$ws.Range("K137").Value = 8639.400000000001
$ws.Range("M137").Value = -3539.400000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5497.7
$ws.Range("I80").Value = 4853
$ws.Range("K80").Value = 4853
$ws.Range("M80").Value = -3855
$ws.Range("H83").Value = 5497.7
$ws.Range("I83").Value = 4853
$ws.Range("K83").Value = 24265
$ws.Range("M83").Value = -19273
$ws.Range("H102").Value = 1466.2307
$ws.Range("I102").Value = 1345.1
$ws.Range("K102").Value = 1345.1
$ws.Range("M102").Value = 276.9000000000001
$ws.Range("H107").Value = 2498.3333
$ws.Range("I107").Value = 2498.3333
$ws.Range("K107").Value = 2498.3333
$ws.Range("M107").Value = -578.3332999999998
$ws.Range("H132").Value = 3193.0715
$ws.Range("I132").Value = 2088.375
$ws.Range("K132").Value = 6265.125
$ws.Range("M132").Value = -3735.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4172.1
$ws.Range("I132").Value = 3996.6667
$ws.Range("J132").Value = 4247.2856
$ws.Range("K132").Value = 11990.0001
$ws.Range("L132").Value = 12741.8568
$ws.Range("M132").Value = -9460.000100000001
$ws.Range("N132").Value = -17801.8568
$ws.Range("H136").Value = 26664664
$ws.Range("I136").Value = 26664664
$ws.Range("K136").Value = 79993992
$ws.Range("M136").Value = -79991442
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2231.84
$ws.Range("I136").Value = 1882.1111
$ws.Range("K136").Value = 5646.3333
$ws.Range("M136").Value = -3096.3333
